$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the lab-section codes in column B (e.g. "CS102L-A1" -> "CS102A-L1")
$ws.Range("B7").Value  = "CS102A-L1"
$ws.Range("B8").Value  = "CS102A-L2"
$ws.Range("B9").Value  = "CS102B-L1"
$ws.Range("B10").Value = "CS102B-L2"
$ws.Range("B13").Value = "CS101A-L1"
$ws.Range("B14").Value = "CS101A-L2"
$ws.Range("B15").Value = "CS101B-L"
$ws.Range("B20").Value = "CS103A-L"
$ws.Range("B21").Value = "CS103B-L"
$ws.Range("B26").Value = "NS104A-L"
$ws.Range("B27").Value = "NS104B-L"
$ws.Range("B30").Value = "CS104A-L"
$ws.Range("B31").Value = "CS104B-L"
$ws.Range("B35").Value = "CS201-L"
$ws.Range("B41").Value = "CS203-L"
$ws.Range("B44").Value = "CS204-L"
$ws.Range("B46").Value = "CS206-L"

# Update the view: scroll position and active cell/selection
$excel.ActiveWindow.ScrollRow = 42
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B62").Select()
